$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.992.24'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.01%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.101.11'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("E4").Value = '  -0.94%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '348.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.79%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5137'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4420'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.98%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '52.32'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08965'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.61%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.167'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '25.28'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.097.52'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.181'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.713'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.18%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '98.74'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.57%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001143'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.73%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.004'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '20.72'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06665'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.23%  '
$ws.Range("E21").Value = '  -0.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.216'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.088.58'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.59'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.62%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.334'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.97%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.344.41'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.75%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.91'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.567'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '162.05'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.17'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.73%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.168'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1057'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.645'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.224'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.969'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.155'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.40%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.09'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02563'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06772'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2270'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.41%  '
$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.46'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.27%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6797'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.43%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.313'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.13'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.67%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6411'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.273'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.71%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000362'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.640'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.54%  '
$ws.Range("E49").Value = '  -2.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '82.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07214'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.63%  '
